$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Portfolio return (column E) input values for each segment row.
$ws.Range("E3").Value = 0.171494307873789
$ws.Range("E4").Value = 0.117184549536677
$ws.Range("E5").Value = 0.0481928182763308
$ws.Range("E6").Value = 0.0942639836674364
$ws.Range("E7").Value = -0.0179878556495654
$ws.Range("E8").Value = 0.0361955361187455
$ws.Range("E9").Value = 0.244557135468765

$ws.Range("E11").Value = 0.0262697558166741
$ws.Range("E12").Value = 0.0474265740434086
$ws.Range("E13").Value = "N/A"

$ws.Range("E15").Value = -0.0227457684867747
$ws.Range("E16").Value = 0.0189999999999999
$ws.Range("E17").Value = "N/A"
$ws.Range("E18").Value = 0.124945601920283
$ws.Range("E19").Value = -0.0198163495055111
$ws.Range("E20").Value = 0.0586768643873218
$ws.Range("E21").Value = -0.00827965837227344
$ws.Range("E22").Value = 0.158484621325823
$ws.Range("E23").Value = 0.0238998922602716

# Update the view: unfreeze scroll position to top-left and move selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E12").Select()

# Restore the window height that the author resized the Excel window to.
$excel.ActiveWindow.Height = 16200
